$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update class diagram text labels (methods/fields) that changed content
# Order chosen to match shared-string table append order from the target workbook
$ws.Range("A7").Value = "_words: List<string>"
$ws.Range("A12").Value = "SplitWords()"
$ws.Range("G7").Value = "_hiddenWords: string"
$ws.Range("G10").Value = "RandomWord()"
$ws.Range("G12").Value = "EndGame()"
$ws.Range("A13").Value = ""

# Copy border/style formatting from G11 (closing box cell) to G12, since the box grew one row
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# G11 is no longer the last row of its box; it should now look like G10 (middle row)
$ws.Range("G10").Copy()
$ws.Range("G11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Move active selection as in diff
$ws.Range("G15").Select()
